# Weekly update: insert a new record at row 26 (most recent price observation),
# pushing all existing rows 26..132 down by one (to 27..133).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 26; Excel shifts rows 26:132 down to 27:133
# and the sheet dimension grows from R132 to R133 automatically.
$ws.Rows("26:26").Insert()

# Populate the newly inserted row 26 with the new weekly data point.
$ws.Cells.Item(26, 1).Value = 10
$ws.Cells.Item(26, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(26, 3).Value = "La Araucanía"
$ws.Cells.Item(26, 4).Value = 45222
$ws.Cells.Item(26, 5).Value = 9
$ws.Cells.Item(26, 6).Value = 300000001
$ws.Cells.Item(26, 7).Value = "Rabanito"
$ws.Cells.Item(26, 8).Value = "Sin especificar"
$ws.Cells.Item(26, 9).Value = "Primera"
$ws.Cells.Item(26, 10).Value = 50
$ws.Cells.Item(26, 11).Value = 9000
$ws.Cells.Item(26, 12).Value = 9000
$ws.Cells.Item(26, 13).Value = 9000
$ws.Cells.Item(26, 14).Value = "`$/docena de paquetes"
$ws.Cells.Item(26, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(26, 16).Value = 750
$ws.Cells.Item(26, 17).Value = 12
$ws.Cells.Item(26, 18).Value = "Hortaliza"
